{"js": "// Add a new bullet to the end of the \"Cose da fare se hai tempo\" list:\n// \"Inserisci nella parte di implementazione delle politiche skip e stop\n// un bel diagramma di sequenza\" \u2014 as a new list item right after the\n// last existing paragraph of the document, sharing its paragraph style\n// and its numbering list.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The paragraph to insert after is the last paragraph in the document\n// (the \"Arricchisci...\" bullet).\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nlastParagraph.load(\"style,isListItem\");\nconst list = lastParagraph.list;\nlist.load(\"id\");\nawait context.sync();\n\nconst newParagraph = lastParagraph.insertParagraph(\n  \"Inserisci nella parte di implementazione delle politiche skip e stop un bel diagramma di sequenza\",\n  Word.InsertLocation.after\n);\n\n// Match the existing bullet's paragraph style (\"Paragrafoelenco\") ...\nnewParagraph.style = lastParagraph.style;\n// ... and keep it in the very same bulleted list (numId) at the same level.\nif (lastParagraph.isListItem) {\n  newParagraph.attachToList(list.id, 0);\n}\n\nawait context.sync();\n", "ps1": "# Add a new bullet to the end of the \"Cose da fare se hai tempo\" list:\n# \"Inserisci nella parte di implementazione delle politiche skip e stop\n# un bel diagramma di sequenza\" \u2014 inserted right after the paragraph\n# that talks about the 58 priority levels (the last bullet already in\n# the document), reusing its paragraph style / list numbering.\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph (the last existing bullet) by its text.\n$searchRange = $d.Content\n$searchRange.Find.Execute(\"Arricchisci la sezione di descrizione\") | Out-Null\n$anchorParagraph = $searchRange.Paragraphs(1)\n$anchorRange = $anchorParagraph.Range\n\n# Insert a new paragraph right after it; Word carries over the paragraph\n# style and list/numbering (numId/ilvl) from the anchor paragraph.\n$anchorRange.InsertParagraphAfter()\n\n# Fill in the text of the freshly inserted (now last) paragraph.\n$newParagraph = $d.Paragraphs.Last\n$newParagraph.Range.Text = \"Inserisci nella parte di implementazione delle politiche skip e stop un bel diagramma di sequenza\"\n"}
